# Edit: "Update The Little Green Web App.pptx"
#  1. Turn the placeholder "#Insert link" quote on slide 7 into a real
#     hyperlinked URL (split across 3 runs, matching PowerPoint's own
#     run-splitting behaviour for mixed/hyperlinked text).
#  2. Delete the trailing "Bye" slide (slide 9).

$p = $ppt.ActivePresentation

# --- 1. Slide 7: replace the quoted placeholder with the live link ---
$s7 = $p.Slides.Item(7)
$sh = $s7.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

$url = "https://jweir211.github.io/LittleGreenWebApp"
$openQuote = [char]0x201C
$closeQuote = [char]0x201D

$part1Text = $openQuote + "https://jweir211.github.io/"
$part2Text = "LittleGreenWebApp"
$part3Text = $closeQuote

$tr.Text = $part1Text + $part2Text + $part3Text

$part1 = $tr.Characters(1, $part1Text.Length)
$part2 = $tr.Characters($part1Text.Length + 1, $part2Text.Length)
$part3 = $tr.Characters($part1Text.Length + $part2Text.Length + 1, $part3Text.Length)

$part1.ActionSettings(1).Hyperlink.Address = $url
$part2.ActionSettings(1).Hyperlink.Address = $url
$part3.ActionSettings(1).Hyperlink.Address = $url

# --- 2. Remove the final "Bye" slide ---
$p.Slides.Item($p.Slides.Count).Delete()
